# Generate Report for Handoff
# Rotate the localization-status report onto a fresh handoff run: the
# generated e2e markdown/xliff artifact names pick up a new GUID + content
# hash, and the handoff timestamps advance a few seconds.

$wb = $excel.ActiveWorkbook

$oldGuid = "f1b6f8a0-260f-4b4d-a252-6981698cd94d"
$newGuid = "ae134d18-3dda-42d6-a21a-1ee067d5c0e3"
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8cc07ef08e39554432ad267c9649312c2f79157/e2e/$oldGuid.md"

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-24 04:55:10"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md") | Out-Null

# ---- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.f11f31b3b97422a3ea8d9b73dc1178f22f74256d.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 04:54:59"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md") | Out-Null

# ---- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.f11f31b3b97422a3ea8d9b73dc1178f22f74256d.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 04:55:10"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md") | Out-Null

Write-Output "Handoff report regenerated."
